$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns (AL1 / AM1): researchAdvisor, otherResearchAdvisor
$ws.Range("AL1").Value = "researchAdvisor"
$ws.Range("AM1").Value = "otherResearchAdvisor"

# Give the new columns custom widths (closest achievable to the source widths)
$ws.Range("AL1").EntireColumn.ColumnWidth = 18.0
$ws.Range("AM1").EntireColumn.ColumnWidth = 17.6667

# Populate row 2 data for the new "researchAdvisor" column, matching the
# existing "otherAdvisor" value in AJ2
$ws.Range("AL2").Value = $ws.Range("AJ2").Value2

# Move/update the active selection to the newly added last header cell
$ws.Range("AM1").Select()
